# poisson_naive versao media ponderada
# Update column A (row index / weighting values) with the new weighted values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2  = 338
    3  = 341
    4  = 343
    5  = 345
    6  = 347
    7  = 350
    8  = 353
    9  = 355
    10 = 357
    11 = 358
    12 = 360
    13 = 361
    14 = 363
    15 = 4
    16 = 80
    17 = 110
    18 = 152
    19 = 193
    20 = 221
    21 = 258
    22 = 303
    23 = 331
    24 = 374
    25 = 405
    26 = 443
    27 = 502
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 1).Value = $values[$row]
}
